$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queries")

# D2 (Start Date) gets the new date string "14-04-2020 00:00:00"
# Leading apostrophe preserves the quote-prefixed text style of the cell
$ws.Range("D2").Value = "'14-04-2020 00:00:00"

# E2 (End Date) keeps the date string that used to be in D2: "15-04-2020 00:00:00"
$ws.Range("E2").Value = "'15-04-2020 00:00:00"

# F2 (Query) gets the updated SQL query text
$query = @'
SELECT M.[AgentID] as [Agent ID],A.[AgentName] as [Agent Name],A.[SupervisorName] as [Supervisor Name],[Channel],[Direction],[DNIS] as [Local Party],
[Ani] as [Remote Party],FORMAT([dbo].[VARCHARTODATETIME](M.ConnectedDateTime),'dd/MM/yyyy HH:mm:ss') as [Interaction Connected Date Time],
FORMAT([dbo].[VARCHARTODATETIME](M.[CreatedDateTime]),'dd/MM/yyyy HH:mm:ss') as [Created Date Time],[SessionID] as [Session ID],CONVERT(varchar, DATEADD(ms, M.QueueTime* 1000, 0), 108) as [Queue Time],
CONVERT(varchar, DATEADD(ms, M.ActiveTime* 1000, 0), 108) as [Active Time],
CONVERT(varchar, DATEADD(ms, M.HoldTime* 1000, 0), 108) as [Hold Time],
CONVERT(varchar, DATEADD(ms, M.ACWTime* 1000, 0), 108) as [ACW Time],
CONVERT(varchar, DATEADD(ms, M.HandleTime* 1000, 0), 108) as [Handle Time],
case when [IsConferenced]=1 Then 'true' when [IsConferenced]=0 then 'false' end as [Is Conferenced],
case when [IsTransfered]=1 Then 'true' when [IsTransfered]=0 then 'false' end as [Is Transfered], 
[TPINTransferReconnected] as [TPIN Transfer Reconnected],[SubChannel] as [Sub Channel],
[SubSessionID] as [Sub Session ID],[InteractionID] as [Interaction ID],[Skill],[SkillName] as [Skill Name],
[DNISName] as [DNIS Name],[TransferedTo] as [Transfered To],[ConferencedTo] as [Conferenced To],
[ConferenceToAgentList] as [Conference To Agent List],[TransferToAgent] as [Transfer To Agent],
[TransferConferenceFromAgent] as [Transfer Conference From Agent] ,[TransferConferenceFromInteraction] as [Transfer Conference From Interaction],
--[OtherData],
FORMAT([dbo].[VARCHARTODATETIME](M.[ClosedDateTime]),'dd/MM/yyyy HH:mm:ss') as [Closed Date Time],
FORMAT([dbo].[VARCHARTODATETIME](M.[DisconnectedDateTime]),'dd/MM/yyyy HH:mm:ss') as [Interaction Disconnected Date Time],[ClosedReason] as [Closed Reason]
--[CIF],[RegisteredMobileNo] as [Registered Mobile No] 
FROM 
(SELECT DISTINCT [User]  AS Ani,
							AgentId as AgentID,
							T.Channel,
							SubChannel,
							T.SessionID AS SessionID,
							SubSessionId as SubSessionID,
							InteractionId as InteractionID,
							T.Direction,
							CreatedDateTime,
							CreatedReason,
							Skill,
							TS.SkillName,
							ISNULL(A.FirstName,'') +' '+ ISNULL(A.LastName,'') AS AgentName,
							Dnis as DNIS,
							DnisName as DNISName,
							IsTransfered,
							IsConferenced,
							IsReconnected AS TPinTransferReconnected,
							IsConferencedTo AS ConferencedTo,
							IsTranferedTo AS TransferedTo,
							CASE WHEN IsTransfered=1 OR IsConferenced=1 THEN  TrasnferConferenceFromAgent ELSE '' END AS TransferConferenceFromAgent,
							CASE WHEN IsTransfered=1 OR IsConferenced=1 THEN  TrasnferConferenceFromInteraction ELSE '' END AS TransferConferenceFromInteraction,
							OtherData,
							ClosedDateTime AS ClosedDateTime,
							ClosedReason,
							CallConnectedTime AS  ConnectedDateTime,
							CallDisconnectedTime AS  DisconnectedDateTime,
							ActiveTime,
							HoldTime,
							TrasnferToAgent AS TransferToAgent,
							ConferenceToAgentList,
							QueueTime,
							AcwTime as ACWTime,
							ActiveTime+HoldTime+AcwTime HandleTime,
							IH.CIF,
							IH.CLID AS RegisteredMobileNo
							FROM TMAC_Interactions T with(nolock)
							INNER JOIN AGT_Agent A with(nolock) ON A.AvayaLoginID=T.AgentId 
							LEFT JOIN AGT_Agent AA with(nolock) ON AA.AvayaLoginID = T.TrasnferConferenceFromAgent
							LEFT JOIN GBL_InteractionHistory IH WITH(NOLOCK) ON IH.SessionID=T.SessionId 
							LEFT JOIN TMAC_Skills TS WITH(NOLOCK) ON TS.SkillExtension=T.Skill
							where 1=1 AND IH.ID IN (SELECT MIN(ID) FROM GBL_InteractionHistory WHERE SESSIONID=T.SessionId AND CLID IS NOT NULL)
							AND  ClosedDateTime>='20200414000000' AND ClosedDateTime<='20200415000000') M
                            INNER JOIN  fn_AgentHierarchy('na','1','1') A  ON A.AgentId=M.AgentID
							order by M.[ClosedDateTime]
'@
$ws.Range("F2").Value = $query

# Restore the row height (Excel caps row height at 409.5pt; the runtime
# auto-fit can overshoot this when the cell gains extra lines of text)
$ws.Rows.Item(2).RowHeight = 409.5

# Update the active selection to F2, matching the edited workbook state
$ws.Range("F2").Select() | Out-Null
